# Generate Report for Handback
# Refresh the handoff/handback timestamps for the file that was just
# handed back (69265fd9-53ea-4c08-8a70-b3598ade048f) in both the
# zh-cn and de-de language sheets, and roll the newest of those two
# timestamps up into the Overview sheet's "Latest HO Xliff Generate
# Date" column for that file's row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-29 18:52:47"
$wsZhCn.Range("K2").Value = "2016-08-29 18:53:14"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-29 18:52:52"
$wsDeDe.Range("K2").Value = "2016-08-29 18:53:21"

# Overview: Latest HO Xliff Generate Date (max across the languages)
$wsOverview.Range("G2").Value = "2016-08-29 18:52:52"
